# "tweaks for the dialogue system and bugfixes"
#
# Adds the next day's log entry to the bottom of the log table on Sheet1:
#   - C5 gets the missing "hours" value (7) for the existing "Movement and
#     dialogue system" row.
#   - A new row 6 is started with the next date (2024-11-22 / serial 45618),
#     formatted the same way as the other date cells in column A.
#   - The active selection is moved to B6 (where the next "Tasks done"
#     entry will be typed in), matching where the author left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing hours figure for the last existing row.
$ws.Range("C5").Value = 7

# Start the next day's row with its date, using the same date format
# ("d-mmm", same style as A2:A5) as the rest of column A.
$ws.Range("A6").Value = 45618
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

# Move the selection to B6, ready for the next "Tasks done" entry.
$ws.Range("B6").Select()
